# feat: add 2022-Q1 data
#
# The existing "总计" (totals) sheet becomes the new "2022-Q1" sheet (per-fund
# holdings for the quarter), and a brand-new "总计" sheet is appended that
# aggregates all quarters including the new 2022-Q1 entry.
#
# NOTE: all Copy()/PasteSpecial() formatting operations are performed BEFORE
# the worksheet is repositioned with Move(); doing so afterwards causes the
# moved sheet (and its move anchor) to lose their styles when used as a copy
# *source* in this runtime, so we carefully order the steps to avoid that.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: repurpose the old "总计" sheet as the new "2022-Q1" sheet, and
# create a fresh blank "总计" sheet.
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Name = "2022-Q1"

$newTotal = $wb.Worksheets.Add()
$newTotal.Name = "总计"

$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$tot = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# Step 2: copy over cell formatting (style "s=2") while sheets are still
# in their just-created positions, before any Move() happens.
# ---------------------------------------------------------------------

# "2022-Q1" sheet: extend the header row's style from D1 into E1:H1
$q1Sheet.Range("D1").Copy()
$q1Sheet.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "总计" sheet: header row style (from "2022-Q1" B1, which still carries the
# original header formatting) and index-column style (from "2022-Q1" A2)
$q1Sheet.Range("B1").Copy()
$tot.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q1Sheet.Range("A2").Copy()
$tot.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 3: reposition the new "总计" sheet right after "2022-Q1".
# ---------------------------------------------------------------------
$tot.Move($null, $q1Sheet)

# Re-fetch fresh references (defensive: avoids any stale-handle issues).
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$tot = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# Step 4: populate the "2022-Q1" sheet (per-fund holdings table)
# ---------------------------------------------------------------------

# Drop the old rows 4:5 (only two data rows are needed now)
$q1Sheet.Range("A4:D5").Clear()

# Header row text
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Row 2: fund 010764
$q1Sheet.Range("A2").Value = 0
$q1Sheet.Range("B2").NumberFormat = "@"
$q1Sheet.Range("B2").Value = "010764"
$q1Sheet.Range("B2").Style = "Normal"
$q1Sheet.Range("C2").NumberFormat = "@"
$q1Sheet.Range("C2").Value = "九泰锐升18个月封闭运作混合"
$q1Sheet.Range("C2").Style = "Normal"
$q1Sheet.Range("D2").NumberFormat = "@"
$q1Sheet.Range("D2").Value = "3.15"
$q1Sheet.Range("D2").Style = "Normal"
$q1Sheet.Range("E2").NumberFormat = "@"
$q1Sheet.Range("E2").Value = "78.81"
$q1Sheet.Range("E2").Style = "Normal"
$q1Sheet.Range("F2").NumberFormat = "@"
$q1Sheet.Range("F2").Value = "5.42"
$q1Sheet.Range("F2").Style = "Normal"
$q1Sheet.Range("G2").NumberFormat = "@"
$q1Sheet.Range("G2").Value = "0.1707"
$q1Sheet.Range("G2").Style = "Normal"
$q1Sheet.Range("H2").Value = 2

# Row 3: fund 009531
$q1Sheet.Range("A3").Value = 1
$q1Sheet.Range("B3").NumberFormat = "@"
$q1Sheet.Range("B3").Value = "009531"
$q1Sheet.Range("B3").Style = "Normal"
$q1Sheet.Range("C3").NumberFormat = "@"
$q1Sheet.Range("C3").Value = "九泰锐和18个月定期开放混合"
$q1Sheet.Range("C3").Style = "Normal"
$q1Sheet.Range("D3").NumberFormat = "@"
$q1Sheet.Range("D3").Value = "1.97"
$q1Sheet.Range("D3").Style = "Normal"
$q1Sheet.Range("E3").NumberFormat = "@"
$q1Sheet.Range("E3").Value = "70.75"
$q1Sheet.Range("E3").Style = "Normal"
$q1Sheet.Range("F3").NumberFormat = "@"
$q1Sheet.Range("F3").Value = "5.79"
$q1Sheet.Range("F3").Style = "Normal"
$q1Sheet.Range("G3").NumberFormat = "@"
$q1Sheet.Range("G3").Value = "0.1141"
$q1Sheet.Range("G3").Style = "Normal"
$q1Sheet.Range("H3").Value = 2

# ---------------------------------------------------------------------
# Step 5: populate the new "总计" sheet (quarter-over-quarter aggregation)
# ---------------------------------------------------------------------
$tot.Range("B1").Value = "日期"
$tot.Range("C1").Value = "持有数量(只)"
$tot.Range("D1").Value = "持有市值(亿元)"

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 2
$tot.Range("D2").Value = 0.28

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2021-Q4"
$tot.Range("C3").Value = 5
$tot.Range("D3").Value = 1.45

$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2021-Q3"
$tot.Range("C4").Value = 4
$tot.Range("D4").Value = 1.84

$tot.Range("A5").Value = 3
$tot.Range("B5").Value = "2021-Q2"
$tot.Range("C5").Value = 2
$tot.Range("D5").Value = 0.03

$tot.Range("A6").Value = 4
$tot.Range("B6").Value = "2021-Q1"
$tot.Range("C6").Value = 1
$tot.Range("D6").Value = 0.01

Write-Host "2022-Q1 data added"
